# Generate Report for Handback
# Adds a new handed-back file ("60e486bb-feec-408d-8370-cc8aa6db5934") to the
# handback-status workbook: one summary row on the "Overview" sheet, and one
# detail row on each of the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$fileId = "60e486bb-feec-408d-8370-cc8aa6db5934"
$mdName = "$fileId.md"
$statusInSync = "Handed back: in sync with en-US"
$handoffReason = "Include"

# ---------------------------------------------------------------------------
# Sheet "Overview": new row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/e2e/$mdName", "", "", $mdName)
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------------
# Sheet "zh-cn": new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$fileId.a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362.zh-cn.xlf"
$zhHandoffDate = "2016-02-22 04:26:24"
$zhHandbackDate = "2016-02-22 04:27:16"

$wsZh.Range("A4").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/e2e/$mdName", "", "", $mdName)

$wsZh.Range("B4").Value = $statusInSync

$wsZh.Range("C4").Value = $zhXlf
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf)

$wsZh.Range("D4").Value = $zhHandoffDate

$wsZh.Range("E4").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/e2e/$mdName", "", "", $mdName)

$wsZh.Range("F4").Value = $zhXlf
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf)

$wsZh.Range("G4").Value = $zhHandbackDate
$wsZh.Range("H4").Value = $handoffReason

# ---------------------------------------------------------------------------
# Sheet "de-de": new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf = "$fileId.a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362.de-de.xlf"
$deHandoffDate = "2016-02-22 04:26:38"
$deHandbackDate = "2016-02-22 04:27:43"

$wsDe.Range("A4").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/e2e/$mdName", "", "", $mdName)

$wsDe.Range("B4").Value = $statusInSync

$wsDe.Range("C4").Value = $deXlf
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf)

$wsDe.Range("D4").Value = $deHandoffDate

$wsDe.Range("E4").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/e2e/$mdName", "", "", $mdName)

$wsDe.Range("F4").Value = $deXlf
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a0de7e4d3b9ac14e6ca85a295cbe25e5c4146362/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf)

$wsDe.Range("G4").Value = $deHandbackDate
$wsDe.Range("H4").Value = $handoffReason

Write-Output "Added handback row for $fileId to Overview, zh-cn, de-de sheets"
